$wb = $excel.ActiveWorkbook

# Rename the first sheet from "olink" to "nanostring" so the tab name
# matches the template contents (nanostring metadata).
$ws1 = $wb.Worksheets.Item("olink")
$ws1.Name = "nanostring"

# The default column width on every sheet shifted by a hair (8.5703125 ->
# 8.578125 chars) when the workbook was last resaved from Excel; mirror
# that here in case it is honoured by the host.
foreach ($ws in $wb.Worksheets) {
    $ws.StandardWidth = 8.578125
}

# Re-select "nanostring" as the active sheet (previously "raw data" was
# active/tabSelected).
$ws1.Activate()

$wb.Save()
